$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "C2"  = 0.0359
    "B3"  = -0.2169
    "C3"  = -0.0381
    "B4"  = -0.1551
    "C4"  = -0.1533
    "B5"  = -0.0219
    "C5"  = -0.0183
    "B6"  = 0.1006
    "C6"  = 0.098
    "B7"  = 0.0905
    "C7"  = 0.0868
    "B8"  = 0.1295
    "C8"  = 0.1286
    "B9"  = -0.0159
    "C9"  = -0.0117
    "B10" = -0.0146
    "C10" = -0.0129
    "B11" = 0.1595
    "C11" = 0.1599
    "B12" = 0.3343
    "C12" = 0.3326
    "B13" = 0.1032
    "C13" = 0.1035
    "C14" = 0.0602
    "C15" = -0.183
    "C17" = -0.3927
    "C18" = 0.0034
    "C20" = 0.2432
    "B24" = -0.1734
    "B27" = 0.6948
    "B30" = -0.0382
    "B32" = -0.262
    "B33" = -0.3126
    "B35" = 0.1324
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
